$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "_old" / "_new" header columns to "_FV2304" / "_FV2310" ---
# Columns A1:J1 carry the "_old" suffix, columns L1:U1 (skipping K1 = "diff")
# carry the "_new" suffix.
$fieldNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($fieldNames[$i])_FV2304"
}

for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($fieldNames[$i])_FV2310"
}

# --- Turn the populated range into an Excel Table (ListObject) ---
$dataRange = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row (split/freeze at row 2) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
